$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 2.680851666666667
$ws.Range("N2").Value = 8.042555
$ws.Range("O2").Value = 0.1074910720871699
$ws.Range("P2").Value = 0.1074910720871699
$ws.Range("Q2").Value = 4.947642218947778
$ws.Range("R2").Value = 44.52877997053
$ws.Range("S2").Value = 0.1074910720871699
$ws.Range("T2").Value = 0.1074910720871699

# Row 3
$ws.Range("M3").Value = 18.51427066666667
$ws.Range("O3").Value = 0.7423457357290222
$ws.Range("P3").Value = 0.7423457357290222
$ws.Range("Q3").Value = 34.16898754317245
$ws.Range("S3").Value = 0.7423457357290222
$ws.Range("T3").Value = 0.7423457357290222

# Row 4
$ws.Range("M4").Value = 3.745104
$ws.Range("N4").Value = 11.235312
$ws.Range("O4").Value = 0.1501631921838079
$ws.Range("P4").Value = 0.1501631921838079
$ws.Range("Q4").Value = 6.911771693728
$ws.Range("R4").Value = 62.205945243552
$ws.Range("S4").Value = 0.1501631921838079
$ws.Range("T4").Value = 0.1501631921838079
